# Removed data from the excel sheet
# - Rename Table1 columns: "Number" -> "PhoneNumber", "Text" -> "TextMessage"
# - Clear the two data rows (row 2 and row 3) on the "Data" sheet, keeping
#   the formatting (style) on column A but dropping all values/strings.
# - Move the active selection on the Data sheet back to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rename the table (ListObject) header columns by writing the header cells
# directly - this also updates the bound ListColumn names.
$ws.Range("A1").Value = "PhoneNumber"
$ws.Range("B1").Value = "TextMessage"

# Clear out the two data rows, but keep column A's existing cell style
# (quote-prefix style used for the phone numbers).
$ws.Range("A2:B3").ClearContents()

# Update the selection to match the new state (B1 instead of B3).
$ws.Activate()
$ws.Range("B1").Select()
